# Adjust responsive design implementation (resolves issue #74)
#
# Appends a new data row (row 69) to each of the four "FE_*" log sheets,
# mirroring the most recent sample row but stamped with the latest capture
# timestamp. Column layout on every sheet:
#   A: time (datetime, formatted as YYYY-MM-DD HH:MM:SS)
#   B: 总长 (hex string)
#   C: ID (hex string)
#   D: 实际长度 (hex string)
#   E: 和校验 (hex string)
#   F: 总长_DEC (number)
#   G: ID_DEC (number)
#   H: 实际长度_DEC (number)
#   I: 和校验_DEC (number)

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"

function Add-LogRow {
    param($ws, $row, $timeValue, $b, $c, $d, $e, $f, $g, $h, $i)

    $ws.Cells.Item($row, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($row, 1).Value2 = $timeValue
    $ws.Cells.Item($row, 2).Value2 = $b
    $ws.Cells.Item($row, 3).Value2 = $c
    $ws.Cells.Item($row, 4).Value2 = $d
    $ws.Cells.Item($row, 5).Value2 = $e
    $ws.Cells.Item($row, 6).Value2 = $f
    $ws.Cells.Item($row, 7).Value2 = $g
    $ws.Cells.Item($row, 8).Value2 = $h
    $ws.Cells.Item($row, 9).Value2 = $i
}

# New timestamp shared by all four sheets' new row (2025-07-17 11:50:22).
$newTime = 45855.49331018519

# The "ID_DEC" (column G) values are huge and need scientific-notation
# parsing; pre-compute them as doubles before calling the function so the
# cast expressions aren't mistaken for a separate statement.
$g1 = [double]"7.598631275147109e+23"
$g2 = [double]"5.68432987514711e+23"
$g3 = [double]"5.68631262647114e+23"
$g4 = [double]"9.85046333984776e+23"

# Sheet 1: FE_LFT_#1
$ws1 = $wb.Worksheets.Item(1)
Add-LogRow $ws1 69 $newTime "0x01,0x7c" "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0," "0x01,0x3C" "0xf" 380 $g1 316 15

# Sheet 2: FE_LFT_#2
$ws2 = $wb.Worksheets.Item(2)
Add-LogRow $ws2 69 $newTime "0x01,0x90" "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78," "0x01,0x4C" "0xe" 400 $g2 332 14

# Sheet 3: FE_PLT_#1
$ws3 = $wb.Worksheets.Item(3)
Add-LogRow $ws3 69 $newTime "0x00,0x6e" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x00,0x63" "0x3" 110 $g3 99 3

# Sheet 4: FE_PLT_#2
$ws4 = $wb.Worksheets.Item(4)
Add-LogRow $ws4 69 $newTime "0x00,0x6e" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x00,0x63" "0x3" 110 $g4 99 3

Write-Output "Added row 69 to all 4 sheets"
